$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4
$ws.Range("A4").Value = "dr. Angyal Gabriella"
$ws.Range("B4").Value = "1036 Budapest,Árpád fejedelem útja 49. I/8."
$ws.Range("C4").Value = "+36 (1) 388 9603, +36 (1) 368 0267, +36 (1) 368 0266"
$ws.Range("D4").Value = "angyalg@mokk.hu"
$ws.Range("E4").Value = "https:\\angyalg.kozjegyzok.mokk.hu"

# Update row 5
$ws.Range("A5").Value = "dr. Steiner Erika Márta"
$ws.Range("B5").Value = "1034 Budapest,Kecske köz 12. "
$ws.Range("C5").Value = "+36 (1) 387 5558, +36 (1) 250 3659"
$ws.Range("D5").Value = "steiner.erika@mokk.hu"
$ws.Range("E5").Value = "https:\\steinererika.kozjegyzok.mokk.hu"

# Append new row 6
$ws.Range("A6").Value = "dr. Kertész Gabriella"
$ws.Range("B6").Value = "1033 Budapest,Miklós utca 11. III/12."
$ws.Range("C6").Value = "+36 (1) 388 9191, +36 (1) 388 2501"
$ws.Range("D6").Value = "kertesz@mokk.hu"
$ws.Range("E6").Value = "https:\\kertesz.kozjegyzok.mokk.hu"

# Append new row 7
$ws.Range("A7").Value = "dr. Barbalics Miklós"
$ws.Range("B7").Value = "1036 Budapest,Árpád fejedelem útja 53/A. I/5."
$ws.Range("C7").Value = "+36 (1) 368 8305, +36 (1) 439 0670"
$ws.Range("D7").Value = "barbalics@mokk.hu"
$ws.Range("E7").Value = "https:\\barbalics.kozjegyzok.mokk.hu"
